# Sample Project / Main.xlsx — "Rules" sheet
# B11 currently holds the text "R40" (rule name). The commit changes it to
# hold the literal text "1" (still a text value, not a number).
#
# A plain  $ws.Range("B11").Value = "1"  would let Excel's type-inference
# treat "1" as a number and store it as <v>1</v> with no shared-string
# reference - not what we want, since the authored change keeps it a text
# cell (t="s") while leaving the cell's existing style untouched.
#
# To force the literal-text interpretation while preserving B11's current
# style, stage the text in a scratch cell formatted as Text ("@"), copy it,
# and use PasteSpecial(xlPasteValues) so only the value/type is transferred
# into B11 - its formatting (style) is left exactly as it was.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$scratch = $ws.Range("ZZ1000")
$scratch.NumberFormat = "@"
$scratch.Value = "1"
$scratch.Copy()

$ws.Range("B11").PasteSpecial(-4163)   # xlPasteValues

$scratch.Clear()
$excel.CutCopyMode = $false
